$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        $systemParts = @()
        $otherParts = @()
        foreach ($p in $parts) {
            if ($p -eq "System" -or $p -eq "system") {
                $systemParts += $p
            } else {
                $otherParts += $p
            }
        }
        $newParts = $systemParts + $otherParts
        $newVal = $newParts -join ", "
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
